$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.496490001678467
$ws.Range("B1").Value = 6.750326633453369
$ws.Range("C1").Value = 5.704201698303223
$ws.Range("D1").Value = 6.721870899200439
$ws.Range("E1").Value = 3.620947122573853
